$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add three new header columns (E:G), copying the header style from D1 ---
$ws.Range("D1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("E1").Value = "original_amount"
$ws.Range("F1").Value = "discount_applied"
$ws.Range("G1").Value = "final_amount"

# --- Fix row 10's phone number: it was stored as text, should be numeric ---
$ws.Range("A10").Value = 76442780

# --- Append the new payment row for 76442780 (Cash) ---
# (leading apostrophe forces text storage, matching the pre-existing
#  "phone stored as text" pattern seen in the other rows of this sheet;
#  reset the style afterwards so the quote-prefix formatting doesn't
#  stick around on the cell)
$ws.Range("A11").Value = "'76442780"
$ws.Range("A11").Style = "Normal"
$ws.Range("C11").Value = "Cash"
$ws.Range("D11").Value = "2025-08-15T10:00:21"
$ws.Range("E11").Value = 200
$ws.Range("F11").Value = 30
$ws.Range("G11").Value = 170
